# added a row for the 2G model
$wb = $excel.ActiveWorkbook
$ws2G = $wb.Worksheets.Item("2G")

# Add the new checklist item ("Only one antenna connector ?") as a new row
# at the bottom of the 2G sheet's checklist.
$ws2G.Range("A15").Value = "Only one antenna connector ?"

# Make the 2G sheet the active tab and move the selection down past the
# newly added row, as it was left after the edit.
$ws2G.Activate()
$ws2G.Range("A16").Select()
